$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("nº1")

# --- Block 1 (rows 3-6) ---
$ws.Range("B3").Value = 1.62
$ws.Range("C3").Value = 0.69
$ws.Range("D3").Value = 1.77
$ws.Range("E3").Value = 0.78

$ws.Range("B4").Value = 0.92
$ws.Range("C4").Value = 732.28
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()

$ws.Range("B5").Value = 0.6
$ws.Range("C5").Value = 0.87
$ws.Range("D5").Value = 0.68
$ws.Range("E5").Value = 1.28

$ws.Range("B6").Value = 0.73
$ws.Range("C6").Value = 0.63
$ws.Range("D6").Value = 0.61
$ws.Range("E6").Value = 0.74

# --- Block 2 (rows 10-13) ---
$ws.Range("B10").Value = 0.8
$ws.Range("C10").Value = 0.72
$ws.Range("D10").Value = 2.38
$ws.Range("E10").Value = 0.92

$ws.Range("B11").Value = 0.83
$ws.Range("C11").Value = 787.74
$ws.Range("D11").ClearContents()
$ws.Range("E11").ClearContents()

$ws.Range("B12").Value = 0.83
$ws.Range("C12").Value = 2.95
$ws.Range("D12").Value = 0.73
$ws.Range("E12").Value = 9.66

$ws.Range("B13").Value = 0.64
$ws.Range("C13").Value = 0.68
$ws.Range("D13").Value = 0.71
$ws.Range("E13").Value = 0.67

# --- Block 3 (rows 17-20) ---
$ws.Range("B17").Value = 1.91
$ws.Range("C17").Value = 0.62
$ws.Range("D17").Value = 1.24
$ws.Range("E17").Value = 1.02

$ws.Range("B18").Value = 0.93
$ws.Range("C18").Value = 787.16
$ws.Range("D18").ClearContents()
$ws.Range("E18").ClearContents()

$ws.Range("B19").Value = 0.67
$ws.Range("C19").Value = 1.67
$ws.Range("D19").Value = 0.67
$ws.Range("E19").Value = 1.4

$ws.Range("B20").Value = 0.73
$ws.Range("C20").Value = 1.21
$ws.Range("D20").Value = 0.65
$ws.Range("E20").Value = 1.93

# --- Block 4 (rows 24-27) ---
$ws.Range("B24").Value = 0.56999999999999995
$ws.Range("C24").Value = 0.69
$ws.Range("D24").Value = 0.56000000000000005
$ws.Range("E24").Value = 0.78

$ws.Range("B25").Value = 0.96
$ws.Range("C25").Value = 972.45
$ws.Range("D25").ClearContents()
$ws.Range("E25").ClearContents()

$ws.Range("B26").Value = 0.69
$ws.Range("C26").Value = 0.66
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 2.2799999999999998

$ws.Range("B27").Value = 0.63
$ws.Range("C27").Value = 0.9
$ws.Range("D27").Value = 0.91
$ws.Range("E27").Value = 1.36

# --- Block 5 (rows 31-34) ---
$ws.Range("B31").Value = 0.99
$ws.Range("C31").Value = 0.82
$ws.Range("D31").Value = 0.67
$ws.Range("E31").Value = 0.75

$ws.Range("B32").Value = 0.95
$ws.Range("C32").Value = 759.38
$ws.Range("D32").ClearContents()
$ws.Range("E32").ClearContents()

$ws.Range("B33").Value = 0.77
$ws.Range("C33").Value = 0.95
$ws.Range("D33").Value = 1.02
$ws.Range("E33").Value = 1.21

$ws.Range("B34").Value = 0.59
$ws.Range("C34").Value = 0.67
$ws.Range("D34").Value = 3.88
$ws.Range("E34").Value = 0.61

# --- Selection ---
$ws.Range("E22").Select()
